$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New override row appended at the bottom of the availability table.
# The "Date" column holds text that looks like an ISO date (e.g. 2026-02-20);
# entering it as a plain string would make Excel auto-convert it to a date
# serial, so we briefly use a leading apostrophe (quote-prefix) to force
# text, then reset the cell style back to Normal so no extra formatting is
# left behind (matching the rest of the sheet, which has no special style on
# these text cells).
$ws.Range("A8").Value = "'2026-02-20"
$ws.Range("A8").Style = "Normal"

$ws.Range("B8").Value = "Limited"
$ws.Range("C8").Value = 12800
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = "Auto-generated from bookings"
